# Edit script: "Added further information to register"
#
# 1) Mark the paragraph "Was ist Clean Code Development?" (and its paragraph
#    mark) as English (en-US) -- adds <w:lang w:val="en-US"/> to both the
#    paragraph-mark run properties (pPr/rPr) and the text run (r/rPr).
#
# 2) Expand the "Hauptteil" section: insert a large block of new
#    sub-headings (Coding Conventions ... Fehlerbehandlung) right after
#    "Die Pfadfinderregel" / "Hinterlasse jeden Code ...", replacing the
#    old, shorter "Coding Conventions ... Errorhandling" block further
#    down (whose last bullet is renamed to "Fehlerbehandlung"). The
#    "_GoBack" bookmark moves along to sit after the new last bullet
#    ("Funktionen sollten klein sein und das Prinzip des SRP befolgen").

$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# (Paragraph 4 = "Was ist Clean Code Development?" -- this InsertXML call
# only rewrites that single paragraph in place, it does not change the
# paragraph count.)
$cleanCodeDevPara = $d.Paragraphs(4)
$cleanCodeDevPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Was ist Clean Code Development?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Change 2 -------------------------------------------------------------
# Find the paragraphs to replace by their text. Paragraph Range.Text carries
# a trailing CR (paragraph mark), so trim it before comparing.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]10)
    if ($t -eq "Hinterlasse jeden Code den du auscheckst besser als du ihn vorgefunden hast.") {
        $startPara = $i
    }
    if ($t -eq "Errorhandling") {
        $endPara = $i
    }
}

$rangeStart = $d.Paragraphs($startPara).Range.Start
$rangeEnd = $d.Paragraphs($endPara).Range.End
$replaceRange = $d.Range($rangeStart, $rangeEnd)
$replaceRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Hinterlasse jeden Code den du auscheckst besser als du ihn vorgefunden hast.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Coding Conventions</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Sind Coding Conventions nötig?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Warum sind sie hilfreich?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wie genau sollte man sie definieren?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Namensgebung</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Variablen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Methoden</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Konstruktoren</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Klassen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Projekte</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Kommentare</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wann sind Kommentare notwendig?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Welche Kommentare sind hilfreich, welche beinhalten nur redundante Informationen?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Formatierung</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Funktionen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Funktionen sollten klein sein und das Prinzip des SRP befolgen</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fehlerbehandlung</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
